# "Create a new 'Net exports' indicator and update the plots"
# This updates the regional summary-statistics table: two regions (the old
# "South Asia" / "Sub-Saharan Africa" rows at the bottom of the table) are
# dropped, "Australia" / "Central Europe and the Baltics" are replaced by
# "East Asia & Pacific" / "Europe & Central Asia", and every region's stats
# are refreshed with newly computed values. The min/max highlight cells
# (cyan = column minimum, green = column maximum) move accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old sheet had 8 data rows (rows 4-11). The new sheet only has 6 data
# rows (rows 4-9), so drop the last two rows.
$ws.Rows(10).Resize(2).Delete()

# Clear the old min/max highlight on E4 (was the old "Australia" row's
# Gross-national-expenditure mean, highlighted cyan as the old minimum) by
# pasting the plain/no-fill format from a never-highlighted neighbour cell.
$ws.Range("C4").Copy()
$ws.Range("E4").PasteSpecial(-4122)   # xlPasteFormats

# Refresh the remaining data rows (4-9) with the new region names + stats.
$data = @(
    @("East Asia & Pacific",        6958.67,  0.25, -1.53, 15128517200000,  0.32, -1.46, 3,    1.32,  1.92),
    @("Europe & Central Asia",     20157.45, -0.47, -1.44, 17540027200000, -0.47, -1.47, 2.89, 1.09,  1.27),
    @("Latin America & Caribbean",  6913.53,  0.04, -1.54,  4067189600000,  0.04, -1.63, 3.88, 0.71,  0.24),
    @("Middle East & North Africa", 5733.24, -0.18, -1.59,  2018788800000,  0.01, -1.73, 3.02, 2.81,  9.55),
    @("North America",             46319.07, -0.08, -1.17, 16274864400000, -0.04, -1.15, 1.98, -0.66, 0.16),
    @("Sub-Saharan Africa",         1259.2,  -0.23, -1.52,  1135243880000, -0.08, -1.66, 5.49, 0.95,  1.25)
)

$row = 4
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $ws.Cells.Item($row, 8).Value = $r[7]
    $ws.Cells.Item($row, 9).Value = $r[8]
    $ws.Cells.Item($row, 10).Value = $r[9]
    $row = $row + 1
}

# Re-apply the min/max highlighting: cyan fill marks the column minimum,
# green fill marks the column maximum (same convention as before the edit).
$cyan = 16776960
$green = 9498256

$ws.Cells.Item(9, 2).Interior.Color = $cyan   # B9  Sub-Saharan Africa (min GDP per capita mean)
$ws.Cells.Item(8, 2).Interior.Color = $green  # B8  North America (max GDP per capita mean)

$ws.Cells.Item(9, 5).Interior.Color = $cyan   # E9  Sub-Saharan Africa (min GNE mean)
$ws.Cells.Item(5, 5).Interior.Color = $green  # E5  Europe & Central Asia (max GNE mean)

$ws.Cells.Item(8, 8).Interior.Color = $cyan   # H8  North America (min inflation mean)
$ws.Cells.Item(9, 8).Interior.Color = $green  # H9  Sub-Saharan Africa (max inflation mean)
